$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (price + volume%) per diff.
$ws.Range("D2").Value = '62.544.61'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.449.04'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'" + '569.80'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").Value = "'" + '145.73'
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'" + '0.527'
$ws.Range("E8").Value = '  -2.08%  '
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").Value = "'" + '0.347'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").Value = "'" + '28.49'
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").Value = "'" + '0.0000173'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").Value = '2.893.19'
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").Value = '62.445.79'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").Value = '2.450.82'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("E18").Value = '  -6.28%  '
$ws.Range("D19").Value = "'" + '10.69'
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").Value = "'" + '320.70'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").Value = "'" + '4.12'
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = "'" + '2.21'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").Value = "'" + '0.999'
$ws.Range("D24").Value = "'" + '9.84'
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").Value = "'" + '64.70'
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("D26").Value = "'" + '642.73'
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").Value = "'" + '0.996'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").Value = '0.0₃0943'
$ws.Range("E29").Value = '  -4.24%  '
$ws.Range("D30").Value = "'" + '1.40'
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("D31").Value = "'" + '7.77'
$ws.Range("E31").Value = '  -3.58%  '
$ws.Range("D32").Value = "'" + '1.80'
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").Value = "'" + '0.999'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -4.27%  '
$ws.Range("D36").Value = "'" + '151.80'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = "'" + '4.60'
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = "'" + '18.47'
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = "'" + '0.362'
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").Value = "'" + '5.27'
$ws.Range("E40").Value = '  -2.70%  '
$ws.Range("D41").Value = "'" + '2.62'
$ws.Range("E41").Value = '  -3.53%  '
$ws.Range("E42").Value = '  -3.75%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '0.0₆0305'
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = "'" + '152.07'
$ws.Range("E45").Value = '  +0.59%  '
$ws.Range("D46").Value = "'" + '15.39'
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("D48").Value = "'" + '0.600'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").Value = "'" + '19.79'
$ws.Range("E49").Value = '  -4.09%  '
$ws.Range("D50").Value = "'" + '0.0500'
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("E51").Value = '  -2.22%  '
